$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying report rows were re-sorted:
#   - row 2 swapped places with row 4 (two "Tretåig hackspett" / Picoides
#     tridactylus records that only differ by id, coordinates, accuracy
#     and the public comment)
#   - row 3 swapped places with row 6 (two lichen/fungus records that
#     only differ by id, taxon id/names and coordinates)
#
# Only the cells that actually differ between the two rows are touched,
# so columns that are identical between the row pairs (and unrelated
# columns such as the date/time text columns) are left completely
# untouched.

# --- swap row 2 <-> row 4 -------------------------------------------------
$a2 = $ws.Range("A2").Value()
$a4 = $ws.Range("A4").Value()
$ws.Range("A2").Value = $a4
$ws.Range("A4").Value = $a2

$q2 = $ws.Range("Q2").Value()
$q4 = $ws.Range("Q4").Value()
$ws.Range("Q2").Value = $q4
$ws.Range("Q4").Value = $q2

$r2 = $ws.Range("R2").Value()
$r4 = $ws.Range("R4").Value()
$ws.Range("R2").Value = $r4
$ws.Range("R4").Value = $r2

$s2 = $ws.Range("S2").Value()
$s4 = $ws.Range("S4").Value()
$ws.Range("S2").Value = $s4
$ws.Range("S4").Value = $s2

# AC2 held a public comment, AC4 was empty - move the comment to AC4
# and clear it out of AC2.
$ac2 = $ws.Range("AC2").Value()
$ws.Range("AC4").Value = $ac2
$ws.Range("AC2").ClearContents()

# --- swap row 3 <-> row 6 -------------------------------------------------
$a3 = $ws.Range("A3").Value()
$a6 = $ws.Range("A6").Value()
$ws.Range("A3").Value = $a6
$ws.Range("A6").Value = $a3

$b3 = $ws.Range("B3").Value()
$b6 = $ws.Range("B6").Value()
$ws.Range("B3").Value = $b6
$ws.Range("B6").Value = $b3

$e3 = $ws.Range("E3").Value()
$e6 = $ws.Range("E6").Value()
$ws.Range("E3").Value = $e6
$ws.Range("E6").Value = $e3

$f3 = $ws.Range("F3").Value()
$f6 = $ws.Range("F6").Value()
$ws.Range("F3").Value = $f6
$ws.Range("F6").Value = $f3

$g3 = $ws.Range("G3").Value()
$g6 = $ws.Range("G6").Value()
$ws.Range("G3").Value = $g6
$ws.Range("G6").Value = $g3

$h3 = $ws.Range("H3").Value()
$h6 = $ws.Range("H6").Value()
$ws.Range("H3").Value = $h6
$ws.Range("H6").Value = $h3

$q3 = $ws.Range("Q3").Value()
$q6 = $ws.Range("Q6").Value()
$ws.Range("Q3").Value = $q6
$ws.Range("Q6").Value = $q3

$r3 = $ws.Range("R3").Value()
$r6 = $ws.Range("R6").Value()
$ws.Range("R3").Value = $r6
$ws.Range("R6").Value = $r3
